# Applies the commit: "Only use today's files and add platform name as a variable."
#  1. Rename the "CJ" sheet to "LOTTE" (only today's / currently-used courier is kept).
#  2. Add a new {PlatformName} variable: used on the LOTTE mapping sheet (G2) and
#     reflected as a renamed header ("Name" -> "Platform Name") on variable_mapping.
#  3. Simplify the LOTTE address mapping to just {long_address}.

$wb = $excel.ActiveWorkbook

$wsMap = $wb.Worksheets.Item("variable_mapping")
$wsCourier = $wb.Worksheets.Item("CJ")

# Rename the courier sheet - it's now the LOTTE template instead of CJ.
$wsCourier.Name = "LOTTE"

# variable_mapping: rename the "Name" header column to "Platform Name".
$wsMap.Range("A1").Value = "Platform Name"

# LOTTE sheet: collapse the address placeholder to a single {long_address} token
# and wire up the new {PlatformName} variable in the previously-empty column G.
$wsCourier.Range("C2").Value = "{long_address}"
$wsCourier.Range("G2").Value = "{PlatformName}"

# Match the refreshed row heights that come from the re-saved template.
$wsCourier.Rows.Item(1).RowHeight = 19.5
$wsCourier.Rows.Item(2).RowHeight = 24.75

# Leave the cursor/active sheet on the LOTTE tab, matching the saved workbook state.
$null = $wsCourier.Range("A4").Select()
